# Weekly update: a new price record (week of 44455) is inserted at row 66,
# pushing the existing historical rows (old 66..96) down by one (new 67..97).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 66; this shifts rows 66-96 down to 67-97,
# and the new dimension will naturally become A1:T97.
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the new weekly record. Columns A,B,C,E,F,G,H,
# I,J,K,L,M,Q,R,T mirror the (unchanged) product/quality metadata that was
# already present on the row being pushed down (now row 67), while D (date),
# N/O/P (min/max/avg price) and S ($/Kg) carry the new observation.
$ws.Range("A66").Value2 = 11
$ws.Range("B66").Value2 = "Vega Monumental Concepción"
$ws.Range("C66").Value2 = "Bíobío"
$ws.Range("D66").Value2 = 44455
$ws.Range("E66").Value2 = 8
$ws.Range("F66").Value2 = "Fruta"
$ws.Range("G66").Value2 = 100108
$ws.Range("H66").Value2 = "Tropicales y subtropicales"
$ws.Range("I66").Value2 = 100108005
$ws.Range("J66").Value2 = "Piña"
$ws.Range("K66").Value2 = "Caramelo"
$ws.Range("L66").Value2 = "Segunda"
$ws.Range("M66").Value2 = 200
$ws.Range("N66").Value2 = 19000
$ws.Range("O66").Value2 = 19500
$ws.Range("P66").Value2 = 19250
$ws.Range("Q66").Value2 = "`$/caja 14 unidades"
$ws.Range("R66").Value2 = "Ecuador"
$ws.Range("S66").Value2 = 1375
$ws.Range("T66").Value2 = 14

# Make sure the new date cell keeps the same date-like numeric formatting as
# the other rows in column D (style index 2 == date custom format).
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
